# AFDP-892 - Store container folder IDs in a single shared table, not in
# columns in the container's table - convert CaseFile module to use
# AcmContainerFolder
#
# The "Assign Alfresco Folder" and "Set Due Date" rules used to read field
# names that lived directly on the case file / container table. Now that
# container folder ids live in the shared AcmContainerFolder table, update
# the rule conditions to match the new object graph.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Assign Alfresco Folder rule: ecmFolderId moved onto containerFolder
$ws.Range("C18").Value = "containerFolder.cmisFolderId == null"

# Set Due Date rule: field name casing corrected to match the bean property
$ws.Range("C21").Value = "dueDate == null"

# Row heights were nudged slightly when the sheet was re-saved
$ws.Rows.Item(18).RowHeight = 13.8
$ws.Rows.Item(21).RowHeight = 13.8

# Restore the view/scroll/selection state that was captured with the file
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.TabRatio = 0.111
$ws.Range("D20").Select()
